$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("F9").Value = 10
$ws.Range("H9").Value = 10

# Row 17
$ws.Range("E17").Value = 112
$ws.Range("F17").Value = 53
$ws.Range("H17").Value = 53

# Row 37
$ws.Range("E37").Value = 51
$ws.Range("F37").Value = 28
$ws.Range("H37").Value = 28

# Row 38
$ws.Range("E38").Value = 71

# Row 41
$ws.Range("E41").Value = 40

# Row 42
$ws.Range("F42").Value = 17
$ws.Range("H42").Value = 17

# Row 45
$ws.Range("F45").Value = 13
$ws.Range("H45").Value = 13

# Row 49
$ws.Range("E49").Value = 64
$ws.Range("F49").Value = 36
$ws.Range("H49").Value = 36

# Row 61
$ws.Range("E61").Value = 28
$ws.Range("F61").Value = 10
$ws.Range("H61").Value = 10

# Row 65
$ws.Range("E65").Value = 28

# Row 71
$ws.Range("F71").Value = 15
$ws.Range("H71").Value = 15

# Row 72
$ws.Range("E72").Value = 40

# Row 74
$ws.Range("E74").Value = 18
$ws.Range("F74").Value = 7
$ws.Range("H74").Value = 7

# Row 75
$ws.Range("F75").Value = 7
$ws.Range("H75").Value = 7

# Row 79
$ws.Range("E79").Value = 37
$ws.Range("F79").Value = 17
$ws.Range("H79").Value = 17
